# Bugfixed evaluation and simulated rt_data for components
# Update ME/MAE/MSE/RMSE/SE and N columns for Q0-Q9 rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = -0.02789355319203728
$ws.Range("C2").Value = 1.215285994576922
$ws.Range("D2").Value = 5.057726240477212
$ws.Range("E2").Value = 2.248938914349879
$ws.Range("F2").Value = 2.270705625358242

# Row 3 (Q1)
$ws.Range("B3").Value = 0.3078495122064144
$ws.Range("C3").Value = 0.9907757216018543
$ws.Range("D3").Value = 3.834933672003624
$ws.Range("E3").Value = 1.958298667722476
$ws.Range("F3").Value = 1.953193691581727

# Row 4 (Q2)
$ws.Range("B4").Value = 0.5026725460285077
$ws.Range("C4").Value = 0.929047527443231
$ws.Range("D4").Value = 2.924062460148327
$ws.Range("E4").Value = 1.709989023399954
$ws.Range("F4").Value = 1.651030220003511

# Row 5 (Q3)
$ws.Range("B5").Value = 0.4047018885367898
$ws.Range("C5").Value = 0.9282943766575162
$ws.Range("D5").Value = 2.415686991994554
$ws.Range("E5").Value = 1.554248047125861
$ws.Range("F5").Value = 1.519276458629088
$ws.Range("G5").Value = 41

# Row 6 (Q4)
$ws.Range("B6").Value = 0.4085208311999072
$ws.Range("C6").Value = 1.124618687279181
$ws.Range("D6").Value = 3.481292747243496
$ws.Range("E6").Value = 1.865822271076079
$ws.Range("F6").Value = 1.850644102371714
$ws.Range("G6").Value = 31

# Row 7 (Q5)
$ws.Range("B7").Value = 0.3717076344828045
$ws.Range("C7").Value = 1.129284152641079
$ws.Range("D7").Value = 3.575005111591263
$ws.Range("E7").Value = 1.890768391842656
$ws.Range("F7").Value = 1.886685660052109
$ws.Range("G7").Value = 29

# Row 8 (Q6)
$ws.Range("B8").Value = 0.3354661335863303
$ws.Range("C8").Value = 1.162998262110728
$ws.Range("D8").Value = 3.668451155248531
$ws.Range("E8").Value = 1.915320118217456
$ws.Range("F8").Value = 1.921634600347395
$ws.Range("G8").Value = 27

# Row 9 (Q7)
$ws.Range("B9").Value = 0.2608331698376827
$ws.Range("C9").Value = 1.384843998384835
$ws.Range("D9").Value = 4.935921458881455
$ws.Range("E9").Value = 2.221693376431918
$ws.Range("F9").Value = 2.266787531232051
$ws.Range("G9").Value = 19

# Row 10 (Q8)
$ws.Range("B10").Value = -0.04940651898553428
$ws.Range("C10").Value = 1.543617672872793
$ws.Range("D10").Value = 7.022163362167988
$ws.Range("E10").Value = 2.649936482666705
$ws.Range("F10").Value = 2.767287288312191
$ws.Range("G10").Value = 12

# Row 11 (Q9)
$ws.Range("B11").Value = -1.179433075504977
$ws.Range("C11").Value = 1.408491936864403
$ws.Range("D11").Value = 4.397925503565186
$ws.Range("E11").Value = 2.097123149356086
$ws.Range("F11").Value = 1.938705471432567
